$d = $word.ActiveDocument

# 1. Title paragraph: merge into one run
$d.Content.Find.Execute(
    "Questions: Rearranging equations involving trigonometry and logarithms",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Questions: Rearranging equations involving trigonometry and logarithms",
    2)

Write-Host "done"
